# Updates cryptocurrency price/volume figures in the "cryptos" worksheet
# to reflect the latest scrape (GitHub Actions run on 2023-01-08 21:39 UTC).
# Values are written as literal text (matching the original inlineStr cells,
# which store numbers/percentages as formatted strings, not numeric values).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $text) {
    $cell = $ws.Range($cellRef)
    # Leading apostrophe forces Excel to store the input as text even
    # though it looks like a number/percentage.
    $cell.Value = "'" + $text
    # Re-apply the Normal style so the quote-prefix formatting Excel
    # added while parsing the text entry does not linger on the cell.
    $cell.Style = "Normal"
}

Set-TextValue "D2" "270.44"
Set-TextValue "E2" "3.37%"
Set-TextValue "E3" "-1.54%"
Set-TextValue "D4" "4.711"
Set-TextValue "E4" "0.03%"
Set-TextValue "D5" "0.06113"
Set-TextValue "E5" "-1.50%"
Set-TextValue "D6" "6.742"
Set-TextValue "E6" "0.25%"
Set-TextValue "D7" "0.8561"
Set-TextValue "E7" "0.72%"
Set-TextValue "D8" "0.8966"
Set-TextValue "E8" "-1.36%"
Set-TextValue "D9" "0.1433"
Set-TextValue "E9" "1.61%"
Set-TextValue "D10" "0.04966"
Set-TextValue "E10" "6.17%"
Set-TextValue "D11" "0.07103"
Set-TextValue "E11" "0.12%"
Set-TextValue "D12" "0.03177"
Set-TextValue "E12" "0.77%"
Set-TextValue "D13" "0.09026"
Set-TextValue "E13" "-0.37%"
Set-TextValue "D14" "0.001542"
Set-TextValue "E14" "0.39%"
Set-TextValue "D15" "0.0006063"
Set-TextValue "E15" "-1.69%"
Set-TextValue "D16" "0.005978"
Set-TextValue "E16" "-2.45%"
Set-TextValue "D17" "3.462"
Set-TextValue "E17" "-0.22%"
Set-TextValue "D18" "3.174"
Set-TextValue "E18" "0.09%"
Set-TextValue "E19" "3.90%"
Set-TextValue "E20" "-0.53%"
Set-TextValue "E21" "-0.68%"
Set-TextValue "D22" "3.836"
Set-TextValue "E22" "-6.23%"
Set-TextValue "D23" "0.04241"
Set-TextValue "E23" "0.23%"
Set-TextValue "D24" "0.001175"
Set-TextValue "E24" "-2.92%"
Set-TextValue "D25" "0.004153"
Set-TextValue "E25" "0.48%"
Set-TextValue "D26" "0.0001200"
Set-TextValue "E26" "-0.08%"
Set-TextValue "D40" "0.03945"
Set-TextValue "E40" "0.98%"
Set-TextValue "D41" "0.1119"
Set-TextValue "E41" "0.46%"
Set-TextValue "D42" "0.004190"
Set-TextValue "E42" "1.36%"
Set-TextValue "D43" "0.002036"
Set-TextValue "E43" "-6.79%"
Set-TextValue "D44" "0.01201"
Set-TextValue "E44" "-13.67%"
Set-TextValue "D45" "0.00005135"
Set-TextValue "E45" "-0.73%"
Set-TextValue "E46" "-0.08%"
Set-TextValue "E47" "-31.82%"
Set-TextValue "D48" "0.9722"
Set-TextValue "E48" "483.35%"
Set-TextValue "E49" "-0.08%"
Set-TextValue "E50" "-0.08%"
